# Scheduled market-data refresh: update crafting profit columns (H:N)
# for the Leve rows whose market prices changed, across the ALC/ARM/BSM/CRP/
# CUL/GSM/LTW/WVR sheets of this Valefor profits workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28: The Writing Is Not on the Wall / Enchanted Silver Ink
$ws.Range("H28").Value = 1491.1333
$ws.Range("I28").Value = 1763
$ws.Range("J28").Value = 1083.3334
$ws.Range("K28").Value = 1763
$ws.Range("L28").Value = 1083.3334
$ws.Range("M28").Value = -1278
$ws.Range("N28").Value = -2053.3334

# Row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 4152
$ws.Range("I76").Value = 3940
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 3940
$ws.Range("L76").Value = 5000
$ws.Range("M76").Value = -3625
$ws.Range("N76").Value = -5630

# Row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 4152
$ws.Range("I79").Value = 3940
$ws.Range("J79").Value = 5000
$ws.Range("K79").Value = 3940
$ws.Range("L79").Value = 5000
$ws.Range("M79").Value = -2848
$ws.Range("N79").Value = -7184

# Row 86: Filling in the Blanks / Enchanted Aurum Regis Ink
$ws.Range("H86").Value = 4108.6665
$ws.Range("I86").Value = 4000
$ws.Range("J86").Value = 4118.5454
$ws.Range("K86").Value = 4000
$ws.Range("L86").Value = 4118.5454
$ws.Range("M86").Value = -2877
$ws.Range("N86").Value = -6364.5454

# Row 89: Ink into Antiquity (L) / Enchanted Aurum Regis Ink
$ws.Range("H89").Value = 4108.6665
$ws.Range("I89").Value = 4000
$ws.Range("J89").Value = 4118.5454
$ws.Range("K89").Value = 20000
$ws.Range("L89").Value = 20592.727
$ws.Range("M89").Value = -14384
$ws.Range("N89").Value = -31824.727

# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 2233831.2
$ws.Range("I132").Value = 2841971.5
$ws.Range("J132").Value = 3983.3333
$ws.Range("K132").Value = 8525914.5
$ws.Range("L132").Value = 11949.9999
$ws.Range("M132").Value = -8523384.5
$ws.Range("N132").Value = -17009.9999

$ws = $wb.Worksheets.Item("ARM")
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 1868.2433
$ws.Range("I132").Value = 1535.04
$ws.Range("J132").Value = 2562.4167
$ws.Range("K132").Value = 4605.12
$ws.Range("L132").Value = 7687.250100000001
$ws.Range("M132").Value = -2075.12
$ws.Range("N132").Value = -12747.2501

$ws = $wb.Worksheets.Item("BSM")
# Row 64: With Bearings Straight / Mythrite Nugget
$ws.Range("H64").Value = 631.6923
$ws.Range("J64").Value = 805.55554
$ws.Range("L64").Value = 805.55554
$ws.Range("N64").Value = -1255.55554

# Row 67: Bearing the Brunt (L) / Mythrite Nugget
$ws.Range("H67").Value = 631.6923
$ws.Range("J67").Value = 805.55554
$ws.Range("L67").Value = 805.55554
$ws.Range("N67").Value = -2365.55554

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 1818.409
$ws.Range("I134").Value = 1063.6364
$ws.Range("J134").Value = 2573.182
$ws.Range("K134").Value = 3190.9092
$ws.Range("L134").Value = 7719.545999999999
$ws.Range("M134").Value = -655.9092000000001
$ws.Range("N134").Value = -12789.546

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 5511.4614
$ws.Range("I16").Value = 3943.625
$ws.Range("J16").Value = 8020
$ws.Range("K16").Value = 3943.625
$ws.Range("L16").Value = 8020
$ws.Range("M16").Value = -3656.625
$ws.Range("N16").Value = -8594

# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 838.7292
$ws.Range("I58").Value = 687.375
$ws.Range("J58").Value = 1141.4375
$ws.Range("K58").Value = 687.375
$ws.Range("L58").Value = 1141.4375
$ws.Range("M58").Value = -484.375
$ws.Range("N58").Value = -1547.4375

# Row 105: Zelkova, My Love / Zelkova Lumber
$ws.Range("H105").Value = 1572.1111
$ws.Range("I105").Value = 941.5
$ws.Range("J105").Value = 2833.3333
$ws.Range("K105").Value = 941.5
$ws.Range("L105").Value = 2833.3333
$ws.Range("M105").Value = 805.5
$ws.Range("N105").Value = -6327.3333

# Row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 5511.4614
$ws.Range("I113").Value = 3943.625
$ws.Range("J113").Value = 8020
$ws.Range("K113").Value = 3943.625
$ws.Range("L113").Value = 8020
$ws.Range("M113").Value = -1773.625
$ws.Range("N113").Value = -12360

# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 1069.5264
$ws.Range("I134").Value = 939.5
$ws.Range("J134").Value = 2174.75
$ws.Range("K134").Value = 2818.5
$ws.Range("L134").Value = 6524.25
$ws.Range("M134").Value = -283.5
$ws.Range("N134").Value = -11594.25

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 838.7292
$ws.Range("I136").Value = 687.375
$ws.Range("J136").Value = 1141.4375
$ws.Range("K136").Value = 2062.125
$ws.Range("L136").Value = 3424.3125
$ws.Range("M136").Value = 487.875
$ws.Range("N136").Value = -8524.3125

$ws = $wb.Worksheets.Item("CUL")
# Row 107: Slippery Service / Frantoio Oil
$ws.Range("H107").Value = 261.47058
$ws.Range("I107").Value = 149.71428
$ws.Range("J107").Value = 339.7
$ws.Range("K107").Value = 449.14284
$ws.Range("L107").Value = 1019.1
$ws.Range("M107").Value = 1470.85716
$ws.Range("N107").Value = -4859.1

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 825.4358999999999
$ws.Range("J131").Value = 1018.2857
$ws.Range("L131").Value = 3054.8571
$ws.Range("N131").Value = -13134.8571

$ws = $wb.Worksheets.Item("GSM")
# Row 99: Needle in a Hingan Stack / Dzo Horn Needle
$ws.Range("H99").Value = 16249.111
$ws.Range("I99").Value = 10206
$ws.Range("J99").Value = 37400
$ws.Range("K99").Value = 10206
$ws.Range("L99").Value = 37400
$ws.Range("M99").Value = -7960
$ws.Range("N99").Value = -41892

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 2584.625
$ws.Range("I102").Value = 1931.7931
$ws.Range("J102").Value = 3581.0527
$ws.Range("K102").Value = 1931.7931
$ws.Range("L102").Value = 3581.0527
$ws.Range("M102").Value = -309.7931000000001
$ws.Range("N102").Value = -6825.0527

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 1906.0714
$ws.Range("I122").Value = 1680.9
$ws.Range("J122").Value = 2110.7727
$ws.Range("K122").Value = 5042.700000000001
$ws.Range("L122").Value = 6332.3181
$ws.Range("M122").Value = -2592.700000000001
$ws.Range("N122").Value = -11232.3181

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 1974.7307
$ws.Range("I7").Value = 2007.7142
$ws.Range("J7").Value = 1936.25
$ws.Range("K7").Value = 2007.7142
$ws.Range("L7").Value = 1936.25
$ws.Range("M7").Value = -1895.7142
$ws.Range("N7").Value = -2160.25

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 3786.9333
$ws.Range("I122").Value = 3708
$ws.Range("K122").Value = 11124
$ws.Range("M122").Value = -8674

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 1974.7307
$ws.Range("I126").Value = 2007.7142
$ws.Range("J126").Value = 1936.25
$ws.Range("K126").Value = 6023.142599999999
$ws.Range("L126").Value = 5808.75
$ws.Range("M126").Value = -3553.142599999999
$ws.Range("N126").Value = -10748.75

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 3130.5264
$ws.Range("I132").Value = 1953.7778
$ws.Range("J132").Value = 4189.6
$ws.Range("K132").Value = 5861.3334
$ws.Range("L132").Value = 12568.8
$ws.Range("M132").Value = -3331.3334
$ws.Range("N132").Value = -17628.8

$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value = 5920.1577
$ws.Range("I107").Value = 967
$ws.Range("J107").Value = 7689.143
$ws.Range("K107").Value = 2901
$ws.Range("L107").Value = 23067.429
$ws.Range("M107").Value = -981
$ws.Range("N107").Value = -26907.429
